# Agile.xlsx - "Maj Docs + Diapo" commit replication
#
# Semantic edit: on sheet "Sprint2", the task row "#8 - Integrer la nouvelle
# implementation du joueur" (old row 12) was removed; the rows below it
# shifted up by one (their text/data moved, but the "#N" labels in column B
# stayed put since they are per-row literal labels). A new tracking day
# (column K, 07/06/2013) was appended to the burndown table, together with
# one more row in the Reste-a-faire% mini table. The named Print_Area for
# Sprint2 shrank by one row to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2")

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# --- New column K: burndown data point for 07/06/2013 ---------------------
Copy-Format "J4" "K4"
$ws.Range("K4").Value = 41432

Copy-Format "J5" "K5"
$ws.Range("K5").Value = 0
Copy-Format "J6" "K6"
$ws.Range("K6").Value = 0
Copy-Format "J7" "K7"
$ws.Range("K7").Value = 0
Copy-Format "J8" "K8"
$ws.Range("K8").Value = 0
Copy-Format "J9" "K9"
$ws.Range("K9").Value = 5
Copy-Format "J10" "K10"
$ws.Range("K10").Value = 5
Copy-Format "J11" "K11"
$ws.Range("K11").Value = 0
Copy-Format "J12" "K12"
$ws.Range("K12").Value = 0
Copy-Format "J13" "K13"
$ws.Range("K13").Value = 0

# --- Remove task row "#8 Integrer la nouvelle implementation du joueur" ---
# Row 12 (C:K) is overwritten with what used to be row 13's content; row 13
# gets what used to be row 14's content; row 14's old content is cleared.
# Column B ("#8", "#9", ...) is left untouched - those labels stay aligned
# to the row position.
$ws.Range("C12").Value = "Création du menu"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 8
$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 4

$ws.Range("C13").Value = "Renseigner le joueur avec un HUD (type de bloc,forme…)"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 5
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 1

# Old row 14 ("#10 Renseigner...") no longer exists as a task row.
$ws.Range("B14:K14").ClearContents()

# --- Totals row moves from 15 to 14, now summing 5:12 instead of 5:13 -----
$ws.Range("E14").Formula = "=SUM(E5:E12)"
$ws.Range("F14").Formula = "=SUM(F5:F12)"
$ws.Range("G14").Formula = "=SUM(G5:G12)"
$ws.Range("H14").Formula = "=SUM(H5:H12)"
$ws.Range("I14").Formula = "=SUM(I5:I12)"
$ws.Range("J14").Formula = "=SUM(J5:J12)"
$ws.Range("K14").Formula = "=SUM(K5:K12)"
$ws.Range("E15:J15").Clear()

# --- "Reste a faire %" mini table: header moves up to row 16, one more ----
# --- data row (22) is added for the new K column ---------------------------
Copy-Format "I17" "I16"
Copy-Format "J17" "J16"
$ws.Range("I16").Value = "Reste à faire %"
$ws.Range("J16").Value = "SPRINT 2"

Copy-Format "I18" "I17"
Copy-Format "J18" "J17"
$ws.Range("I17").Formula = "=F14/`$E`$14"
$ws.Range("J17").Formula = "=F`$4"

Copy-Format "I19" "I18"
Copy-Format "J19" "J18"
$ws.Range("I18").Formula = "=G14/`$E`$14"
$ws.Range("J18").Formula = "=G4"

Copy-Format "I20" "I19"
Copy-Format "J20" "J19"
$ws.Range("I19").Formula = "=H14/`$E`$14"
$ws.Range("J19").Value = 41428

Copy-Format "I21" "I20"
Copy-Format "J21" "J20"
$ws.Range("I20").Formula = "=I14/`$E`$14"
$ws.Range("J20").Value = 41429

Copy-Format "I22" "I21"
Copy-Format "J22" "J21"
$ws.Range("I21").Formula = "=J14/`$E`$14"
$ws.Range("J21").Value = 41430

Copy-Format "I21" "I22"
Copy-Format "J21" "J22"
$ws.Range("I22").Formula = "=K14/`$E`$14"
$ws.Range("J22").Value = 41432

# --- Print area: Sprint2 used to end at row 35, now ends at row 34 --------
$ws.PageSetup.PrintArea = "B2:K34"

# --- Selection / scroll position mirrors the committed view ---------------
$ws.Application.Goto($ws.Range("A10"), $false)
$ws.Range("I26").Select() | Out-Null
